$wb = $excel.ActiveWorkbook

# 1) Sales sheet: move selection from E5 to C5
$wb.Worksheets.Item("Sales").Activate()
$wb.Worksheets.Item("Sales").Range("C5").Select()

# 2) Reorder sheets: move "Returns" to sit after "PG2" (i.e. before "Invoices to collect").
#    NOTE: worksheet references are index-bound, so re-fetch everything by name
#    after this call - stale handles now point at whatever sheet occupies their
#    old slot.
$wb.Worksheets.Item("Returns").Move($null, $wb.Worksheets.Item("PG2"))

# 3) Returns sheet: selection F2 -> F4, column F width -> 32 (raw xlsx width)
$wb.Worksheets.Item("Returns").Activate()
$wb.Worksheets.Item("Returns").Columns.Item(6).ColumnWidth = 31.16666666666667
$wb.Worksheets.Item("Returns").Range("F4").Select()

# 4) Invoices to collect: keep selection at A4 (it loses tabSelected once another
#    sheet becomes active below).
$wb.Worksheets.Item("Invoices to collect").Activate()
$wb.Worksheets.Item("Invoices to collect").Range("A4").Select()

# 5) Payment methods: becomes the active tab, selection F6 -> F5, and F4 value
#    "07 August 2018" -> "27 August 2018" (leading apostrophe keeps it literal
#    text like the original cell, preserving its quotePrefix style).
$wb.Worksheets.Item("Payment methods").Activate()
$wb.Worksheets.Item("Payment methods").Range("F4").Value = "'27 August 2018"
$wb.Worksheets.Item("Payment methods").Range("F5").Select()
